# Weekly driver report update for 2025-04-28
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Intel(R) Wi-Fi 6E AX211 160MHz - 23.90.0.2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 102
$ws.Range("D3").Value = 98.59999999999999

# Row 4: Totals
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 102

# Row 12: Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1
$ws.Range("B12").Value = 11140

# Row 13: Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3
$ws.Range("B13").Value = 14487
